# Fill in the computed Retention ratio and Answer Recall Average (ARA)
# metric values in the summary table. Each target row has its value in
# the last cell of the row, which starts out empty.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of "label text contained in row's first cell" -> value to write
# into that row's last (score) cell.
$targets = @(
    @{ Label = "Ratio"; Value = "0.5" },
    @{ Label = "Answer Recall Lenient (ARL)"; Value = "0.375" },
    @{ Label = "Answer Recall Strict (ARS)"; Value = "0.25" },
    @{ Label = "Answer Recall Average (ARA)"; Value = "0.3125" }
)

foreach ($target in $targets) {
    for ($i = 1; $i -le $t.Rows.Count; $i++) {
        $row = $t.Rows.Item($i)
        $labelCell = $row.Cells.Item(1)
        # Cell text ends with a paragraph mark (CR) + cell mark (BEL);
        # strip those control characters before comparing.
        $labelText = $labelCell.Range.Text.TrimEnd([char]13, [char]7).Trim()
        if ($labelText -eq $target.Label) {
            $scoreCell = $row.Cells.Item($row.Cells.Count)
            $r = $scoreCell.Range
            $r.Text = $target.Value
            $r2 = $scoreCell.Range
            $r2.Font.Bold = $true
            $r2.Font.Size = 12
            $r2.Font.SizeBi = 12
            break
        }
    }
}
